$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet (tab name changes from "Quantidade de Alunos" to "QuantidadeAlunos")
$ws.Name = "QuantidadeAlunos"

# Re-label and reorder the header row values:
#   A1 -> "SchoolUnit"
#   B1 -> "NumberOfStudents"
#   C1 -> "Code"
$ws.Range("A1").Value = "SchoolUnit"
$ws.Range("B1").Value = "NumberOfStudents"
$ws.Range("C1").Value = "Code"

# Update the active-cell selection to H12
$ws.Range("H12").Select()
